$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF), matching the existing
# header formatting (bold, centered, bordered) used by A1:H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new I/J columns, rows 2-70.
$iVals = @(9,7,9,6,6,6,8,9,6,8,6,6,6,7,7,6,8,8,8,6,8,6,6,8,7,6,7,8,10,7,8,8,7,6,8,8,7,7,7,7,7,8,6,6,8,7,9,8,6,7,7,10,7,8,8,7,7,8,6,6,7,5,7,6,5,8,7,8,5)
$jVals = @(9,7,9,6,7,6,8,9,6,8,7,6,7,7,7,6,8,9,8,7,8,7,7,9,7,7,7,8,10,7,8,8,7,6,9,8,7,7,7,7,7,8,7,6,8,7,9,8,7,7,7,10,7,8,8,7,7,8,7,7,8,7,7,7,6,8,7,8,5)

for ($r = 2; $r -le 70; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
    $ws.Cells.Item($r, 10).Value = $jVals[$idx]
}
